$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose text could be mis-parsed as numbers by Excel need a Text
# number format applied first, then cleared again afterwards so the
# cell keeps its original (default) style but the value stays a string.
function Set-TextValue($cell, $value) {
    $r = $ws.Range($cell)
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.ClearFormats()
}

# Row 2
Set-TextValue 'D2' '67.890.55'
Set-TextValue 'E2' '  -2.24%  '

# Row 3
Set-TextValue 'D3' '3.802.25'
Set-TextValue 'E3' '  +1.02%  '

# Row 4
Set-TextValue 'E4' '  +0.04%  '

# Row 5
Set-TextValue 'D5' '601.45'
Set-TextValue 'E5' '  -2.17%  '

# Row 6
Set-TextValue 'D6' '172.45'
Set-TextValue 'E6' '  -3.60%  '

# Row 7
Set-TextValue 'D7' '3.799.42'
Set-TextValue 'E7' '  +1.04%  '

# Row 8
Set-TextValue 'E8' '  +0.01%  '

# Row 9
Set-TextValue 'D9' '0.529'
Set-TextValue 'E9' '  +0.25%  '

# Row 10
Set-TextValue 'E10' '  -4.57%  '

# Row 11
Set-TextValue 'E11' '  -5.79%  '

# Row 12
Set-TextValue 'D12' '0.467'
Set-TextValue 'E12' '  -3.68%  '

# Row 13
Set-TextValue 'D13' '38.77'
Set-TextValue 'E13' '  -3.55%  '

# Row 14
Set-TextValue 'D14' '0.0000244'
Set-TextValue 'E14' '  -3.53%  '

# Row 15
Set-TextValue 'D15' '4.443.42'
Set-TextValue 'E15' '  +1.20%  '

# Row 16
Set-TextValue 'D16' '3.805.63'
Set-TextValue 'E16' '  +1.13%  '

# Row 17
Set-TextValue 'D17' '67.859.39'
Set-TextValue 'E17' '  -2.37%  '

# Row 18
Set-TextValue 'D18' '7.26'
Set-TextValue 'E18' '  -3.79%  '

# Row 19
Set-TextValue 'E19' '  -3.93%  '

# Row 20
Set-TextValue 'D20' '17.25'
Set-TextValue 'E20' '  +5.35%  '

# Row 21
Set-TextValue 'D21' '494.41'
Set-TextValue 'E21' '  -2.93%  '

# Row 22
Set-TextValue 'D22' '9.16'
Set-TextValue 'E22' '  -2.08%  '

# Row 23
Set-TextValue 'D23' '0.741'
Set-TextValue 'E23' '  +1.38%  '

# Row 24
Set-TextValue 'D24' '85.99'
Set-TextValue 'E24' '  -0.62%  '

# Row 25
Set-TextValue 'E25' '  -5.21%  '

# Row 26
Set-TextValue 'D26' '0.0000145'
Set-TextValue 'E26' '  +8.42%  '

# Row 27
Set-TextValue 'D27' '12.35'
Set-TextValue 'E27' '  -4.03%  '

# Row 28
Set-TextValue 'D28' '10.23'
Set-TextValue 'E28' '  -3.79%  '

# Row 29
Set-TextValue 'E29' '  +0.09%  '

# Row 30
Set-TextValue 'D30' '2.97'
Set-TextValue 'E30' '  -0.66%  '

# Row 31
Set-TextValue 'D31' '2.44'
Set-TextValue 'E31' '  -3.14%  '

# Row 32
Set-TextValue 'D32' '32.77'
Set-TextValue 'E32' '  +6.97%  '

# Row 33
Set-TextValue 'D33' '7.86'
Set-TextValue 'E33' '  -1.67%  '

# Row 34
Set-TextValue 'E34' '  -4.20%  '

# Row 35
Set-TextValue 'E35' '  +0.07%  '

# Row 36
Set-TextValue 'E36' '  -3.93%  '

# Row 37
Set-TextValue 'D37' '5.83'
Set-TextValue 'E37' '  -5.29%  '

# Row 38
Set-TextValue 'E38' '  -2.87%  '

# Row 39
Set-TextValue 'D39' '462.23'
Set-TextValue 'E39' '  +1.12%  '

# Row 40
Set-TextValue 'E40' '  -5.45%  '

# Row 41
Set-TextValue 'B41' 'Stacks'
Set-TextValue 'C41' 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
Set-TextValue 'D41' '2.02'
Set-TextValue 'E41' '  -2.77%  '

# Row 42
Set-TextValue 'B42' 'OKB'
Set-TextValue 'C42' 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
Set-TextValue 'D42' '49.08'
Set-TextValue 'E42' '  -1.46%  '

# Row 43
Set-TextValue 'D43' '2.85'
Set-TextValue 'E43' '  -4.24%  '

# Row 44
Set-TextValue 'D44' '8.43'
Set-TextValue 'E44' '  -1.74%  '

# Row 45
Set-TextValue 'D45' '41.05'
Set-TextValue 'E45' '  -8.82%  '

# Row 47
Set-TextValue 'D47' '2.848.87'
Set-TextValue 'E47' '  -3.65%  '

# Row 48
Set-TextValue 'D48' '139.93'
Set-TextValue 'E48' '  +0.61%  '

# Row 49
Set-TextValue 'D49' '0.0352'
Set-TextValue 'E49' '  -2.40%  '

# Row 50
Set-TextValue 'D50' '26.01'
Set-TextValue 'E50' '  -5.20%  '

# Row 51
Set-TextValue 'D51' '24.00'
Set-TextValue 'E51' '  +10.79%  '
